$wb = $excel.ActiveWorkbook

# --- Update the "Hoja1" sheet: cell A1 conversion text ---
$ws1 = $wb.Worksheets.Item("Hoja1")

$newText = "Conversión del día 💰`r`n✅ Dólar paralelo: 68`r`n`r`nBinance`r`n✅ 1000 Bs = 2.68 = 10080.43 pesos`r`n✅ 10080.43 pesos = 2.68 = 951.74 Bs`r`n`r`nPromedio competencia`r`n✅ Tasa pesos: 20`r`n✅ Tasa Bs: 20`r`n✅ % Ganancia: 20%"

$ws1.Range("A1").Value = $newText

# --- Update the "tasas" sheet: rate figures in N10, O10, N12, O12 ---
$ws2 = $wb.Worksheets.Item("tasas")

$ws2.Range("N10").Value = 373
$ws2.Range("O10").Value = 3760
$ws2.Range("N12").Value = 3760
$ws2.Range("O12").Value = 355
